$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'54.834.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +8.23%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.409.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +6.88%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'474.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +12.56%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'138.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +20.01%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.12%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.501"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +11.52%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.439.11"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +7.88%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0952"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +12.34%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'5.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +8.17%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.322"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +10.39%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.122"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.19%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.840.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +7.71%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'54.922.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +8.28%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'20.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +11.98%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +17.69%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.442.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +8.30%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +11.01%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'9.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +17.58%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'310.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +7.86%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.991"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.18%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +12.87%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'56.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +8.98%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.49%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.400"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +11.44%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.161"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +20.44%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.554.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +9.54%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +10.40%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0₃0761"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +20.51%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.997"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.01%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'148.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.31%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'17.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +9.62%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +13.83%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.12"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +12.22%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +16.73%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +9.05%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.826"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +11.86%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'33.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +6.35%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.997"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.03%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'Filecoin"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'3.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +12.31%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'Mantle"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.598"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +9.98%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +12.24%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +17.09%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'4.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +25.32%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'10.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.32%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'252.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +34.81%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0891"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +14.53%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0220"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +12.26%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.902.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.40%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'16.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +12.32%  "
$ws.Range("E51").Style = "Normal"
